$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new value mapping derived from the diff (symbol list refresh).
$updates = [ordered]@{
    'D2' = '305.50'
    'E2' = '-0.85%'
    'G2' = '23'
    'D3' = '39.08'
    'E3' = '7.63%'
    'G3' = '23'
    'D4' = '5.114'
    'E4' = '1.12%'
    'G4' = '23'
    'D5' = '0.08069'
    'E5' = '-0.50%'
    'G5' = '23'
    'D6' = '1.935'
    'E6' = '-2.66%'
    'G6' = '23'
    'B7' = 'KuCoinToken'
    'C7' = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
    'D7' = '8.047'
    'E7' = '2.69%'
    'G7' = '23'
    'B8' = 'MXToken'
    'C8' = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    'D8' = '0.9263'
    'E8' = '-0.38%'
    'G8' = '23'
    'B9' = 'LiechtensteinCryptoassetsExchange'
    'C9' = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
    'D9' = '0.1350'
    'E9' = '-8.44%'
    'G9' = '23'
    'B10' = 'WazirX'
    'C10' = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
    'D10' = '0.1913'
    'E10' = '-0.88%'
    'G10' = '23'
    'B11' = 'MandalaExchangeToken'
    'C11' = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
    'D11' = '0.08993'
    'E11' = '-1.17%'
    'G11' = '23'
    'B12' = 'BitrueCoin'
    'C12' = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
    'D12' = '0.03519'
    'E12' = '0.02%'
    'G12' = '23'
    'B13' = 'BitMartToken'
    'C13' = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
    'D13' = '0.09760'
    'E13' = '-0.99%'
    'G13' = '23'
    'B14' = 'BitForexToken'
    'C14' = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
    'D14' = '0.001406'
    'E14' = '-1.00%'
    'G14' = '23'
    'B15' = 'TigerCash'
    'C15' = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
    'D15' = '0.005909'
    'E15' = '-6.35%'
    'G15' = '23'
    'B16' = 'LEO'
    'C16' = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
    'D16' = '3.764'
    'E16' = '-2.22%'
    'G16' = '23'
    'B17' = 'GateToken'
    'C17' = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
    'D17' = '4.208'
    'E17' = '1.22%'
    'G17' = '23'
    'E18' = '-1.47%'
    'G18' = '23'
    'D19' = '0.3461'
    'E19' = '0.37%'
    'G19' = '23'
    'D20' = '0.1319'
    'E20' = '-0.61%'
    'G20' = '23'
    'D21' = '4.682'
    'E21' = '-2.69%'
    'G21' = '23'
    'D22' = '0.2417'
    'E22' = '3.05%'
    'G22' = '23'
    'D23' = '0.04367'
    'E23' = '-0.10%'
    'G23' = '23'
    'D24' = '0.001207'
    'E24' = '-2.39%'
    'G24' = '23'
    'D25' = '0.004268'
    'E25' = '2.62%'
    'G25' = '23'
    'D26' = '0.0001303'
    'E26' = '0.01%'
    'G26' = '23'
    'G27' = '23'
    'G28' = '23'
    'G29' = '23'
    'G30' = '23'
    'G31' = '23'
    'G32' = '23'
    'G33' = '23'
    'G34' = '23'
    'G35' = '23'
    'G36' = '23'
    'G37' = '23'
    'G38' = '23'
    'D39' = '0.02027'
    'E39' = '-1.76%'
    'G39' = '23'
    'D40' = '0.05034'
    'E40' = '-1.25%'
    'G40' = '23'
    'D41' = '0.007522'
    'E41' = '0.87%'
    'G41' = '23'
    'D42' = '0.009666'
    'E42' = '-4.58%'
    'G42' = '23'
    'D43' = '0.1343'
    'E43' = '-1.72%'
    'G43' = '23'
    'D44' = '0.002094'
    'E44' = '-1.41%'
    'G44' = '23'
    'D45' = '0.009787'
    'E45' = '1.01%'
    'G45' = '23'
    'D46' = '0.00006226'
    'E46' = '-0.88%'
    'G46' = '23'
    'D47' = '0.00000000752'
    'E47' = '0.02%'
    'G47' = '23'
    'D48' = '0.002874'
    'G48' = '23'
    'E49' = '12.50%'
    'G49' = '23'
    'D50' = '0.00002104'
    'E50' = '0.02%'
    'G50' = '23'
    'D51' = '0.0002004'
    'E51' = '0.02%'
    'G51' = '23'
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    # Prefix with an apostrophe so Excel stores the value as literal text
    # (preserves formatting such as trailing zeros, leading zeros and the % sign)
    $cell.Value = "'" + $updates[$ref]
    # Drop any quote-prefix / number formatting Excel may have implicitly applied
    $cell.Style = "Normal"
}
